# Add a new "2022-Q3" sheet (with the same layout/formatting as the
# existing "2022-Q1" sheet) right after the "总计" sheet, fill it with
# the 2022-Q3 fund holding data, and update the "总计" summary sheet so
# a new row for 2022-Q3 is inserted at the top of the data (pushing the
# existing rows down by one).

$wb = $excel.ActiveWorkbook

# Writes $value into $cell as literal text (keeps "11.66"-looking values
# from being auto-converted to numbers), without leaving the cell's
# number format permanently changed.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet by duplicating "2022-Q1" (this
#    keeps header/row styles identical) and placing it right after the
#    first sheet ("总计").
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$template = $wb.Worksheets.Item("2022-Q1")
$template.Copy($null, $summary) | Out-Null

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Extend the formatting of the existing data row (row 2) down to rows
# 3-7 so every data row shares the same look (bold/border on column A).
$q3.Range("A2:H2").Copy() | Out-Null
$q3.Range("A3:H7").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# Helper data for the new sheet: each entry is
# (fund code, fund name, fund size, total stock position, position
#  ratio, market value held, position rank)
# ---------------------------------------------------------------------
$q3Rows = @(
    @("007553", "中信建投医改灵活配置混合C", "11.66", "94.99", "3.86", "0.4501", 10),
    @("002408", "中信建投医改灵活配置混合A", "10.95", "94.99", "3.86", "0.4227", 10),
    @("015139", "泰康医疗健康股票A",         "0.34",  "85.50", "3.05", "0.0104", 6),
    @("012432", "国投瑞银安泰混合C",         "1.00",  "32.06", "0.85", "0.0085", 9),
    @("015140", "泰康医疗健康股票C",         "0.13",  "85.50", "3.05", "0.0040", 6),
    @("012431", "国投瑞银安泰混合A",         "0.00",  "32.06", "0.85", $null,    9)
)

for ($i = 0; $i -lt $q3Rows.Count; $i++) {
    $row = 2 + $i
    $item = $q3Rows[$i]

    $q3.Cells.Item($row, 1).Value = $i

    Set-TextValue $q3.Cells.Item($row, 2) $item[0]
    Set-TextValue $q3.Cells.Item($row, 3) $item[1]
    Set-TextValue $q3.Cells.Item($row, 4) $item[2]
    Set-TextValue $q3.Cells.Item($row, 5) $item[3]
    Set-TextValue $q3.Cells.Item($row, 6) $item[4]

    if ($item[5] -eq $null) {
        # last row's market value is a genuine numeric 0, not text "0.00"
        $q3.Cells.Item($row, 7).Value = 0
    } else {
        Set-TextValue $q3.Cells.Item($row, 7) $item[5]
    }

    $q3.Cells.Item($row, 8).Value = $item[6]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert the 2022-Q3 totals as the new
#    first data row and shift the previous rows down by one.
# ---------------------------------------------------------------------
$totalRows = @(
    @("2022-Q3", 6, 0.9),
    @("2022-Q1", 1, 0.05),
    @("2021-Q2", 2, 0.68),
    @("2021-Q1", 7, 0.84),
    @("2020-Q4", 1, 0.31)
)

# The summary sheet only had formatting defined through row 5; extend
# it down to the new row 6 so the A-column keeps its bold/border style.
$summary.Range("A5:D5").Copy() | Out-Null
$summary.Range("A6:D6").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $row = 2 + $i
    $item = $totalRows[$i]

    $summary.Cells.Item($row, 1).Value = $i
    $summary.Cells.Item($row, 2).Value = $item[0]
    $summary.Cells.Item($row, 3).Value = $item[1]
    $summary.Cells.Item($row, 4).Value = $item[2]
}

# Restore "2020-Q4" (the sheet that was active before this edit) as the
# active sheet, matching the original workbook's tab selection.
$wb.Worksheets.Item("2020-Q4").Activate() | Out-Null
